$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.829.48"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "'1.900.35"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'0.7700"
$ws.Range("E5").Value = "  +4.04%  "
$ws.Range("D6").Value = "'240.13"
$ws.Range("E6").Value = "  -1.62%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.3055"
$ws.Range("E8").Value = "  -2.40%  "
$ws.Range("D9").Value = "'25.47"
$ws.Range("E9").Value = "  -5.49%  "
$ws.Range("D10").Value = "'0.06853"
$ws.Range("E10").Value = "  -1.89%  "
$ws.Range("D11").Value = "'0.07983"
$ws.Range("D12").Value = "'1.905.25"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "'0.7374"
$ws.Range("E13").Value = "  -5.59%  "
$ws.Range("D14").Value = "'5.178"
$ws.Range("E14").Value = "  -2.14%  "
$ws.Range("D15").Value = "'91.38"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("D16").Value = "'29.856.99"
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("D17").Value = "'13.76"
$ws.Range("E17").Value = "  -4.58%  "
$ws.Range("D18").Value = "'5.897"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "'246.31"
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("D20").Value = "'0.000007701"
$ws.Range("E20").Value = "  -1.86%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "'2.149.24"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "'6.946"
$ws.Range("E24").Value = "  -3.98%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "'9.270"
$ws.Range("E25").Value = "  -2.00%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'166.57"
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("E27").Value = "  -2.03%  "
$ws.Range("D28").Value = "'0.1286"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "'2.028"
$ws.Range("E29").Value = "  -2.10%  "
$ws.Range("D30").Value = "'1.392"
$ws.Range("E30").Value = "  +2.74%  "
$ws.Range("D31").Value = "'1.510"
$ws.Range("E31").Value = "  -2.49%  "
$ws.Range("D32").Value = "'4.271"
$ws.Range("D34").Value = "'0.05251"
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("E35").Value = "  -4.44%  "
$ws.Range("D36").Value = "'0.7278"
$ws.Range("E36").Value = "  -3.17%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "'0.01908"
$ws.Range("E38").Value = "  -2.02%  "
$ws.Range("D39").Value = "'2.779"
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("D40").Value = "'6.188"
$ws.Range("E40").Value = "  -3.02%  "
$ws.Range("D41").Value = "'0.4412"
$ws.Range("E41").Value = "  -2.50%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "'0.8365"
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("D45").Value = "'1.877"
$ws.Range("E45").Value = "  -4.81%  "
$ws.Range("D46").Value = "'7.598"
$ws.Range("E46").Value = "  -3.83%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'9.811"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'100.29"
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("D49").Value = "'2.054.55"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").Value = "'36.14"
$ws.Range("E50").Value = "  -2.90%  "
$ws.Range("D51").Value = "'916.21"
$ws.Range("E51").Value = "  -2.51%  "
